$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date values for rows 2-42 advance by one day:
# 45714 -> 45715 (2025-02-26 -> 2025-02-27)
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 3).Value = 45715
}
